$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.407.04'
$ws.Range('E2').Value = '  -0.43%  '
$ws.Range('D3').Value = '1.724.69'
$ws.Range('E3').Value = '  -0.18%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '243.70'
$ws.Range('E5').Value = '  -0.26%  '
$ws.Range('E6').Value = '  +0.06%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4897'
$ws.Range('E7').Value = '  +1.85%  '
$ws.Range('E8').Value = '  -2.29%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06193'
$ws.Range('E9').Value = '  +0.12%  '
$ws.Range('D10').Value = '1.720.54'
$ws.Range('E10').Value = '  -0.82%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07004'
$ws.Range('E11').Value = '  -2.51%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.50'
$ws.Range('E12').Value = '  -0.57%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.572'
$ws.Range('E13').Value = '  +0.94%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5992'
$ws.Range('E14').Value = '  -2.00%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '77.19'
$ws.Range('E15').Value = '  +0.04%  '
$ws.Range('E16').Value = '  +0.07%  '
$ws.Range('D17').Value = '26.417.24'
$ws.Range('E17').Value = '  -0.42%  '
$ws.Range('E18').Value = '  +0.07%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007114'
$ws.Range('E19').Value = '  +2.58%  '
$ws.Range('E20').Value = '  -1.70%  '
$ws.Range('D21').Value = '1.943.07'
$ws.Range('E21').Value = '  -0.68%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.465'
$ws.Range('E22').Value = '  -1.36%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.587'
$ws.Range('E23').Value = '  -2.38%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.153'
$ws.Range('E24').Value = '  -1.88%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '137.36'
$ws.Range('E25').Value = '  +0.14%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.21'
$ws.Range('E26').Value = '  -0.77%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.392'
$ws.Range('E27').Value = '  -0.36%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '106.99'
$ws.Range('E28').Value = '  -0.28%  '
$ws.Range('E29').Value = '  -4.49%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.942'
$ws.Range('E30').Value = '  -0.49%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.07934'
$ws.Range('E31').Value = '  -1.09%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.675'
$ws.Range('E32').Value = '  -0.48%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04532'
$ws.Range('E33').Value = '  +0.29%  '
$ws.Range('B34').Value = 'Frax'
$ws.Range('C34').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.9995'
$ws.Range('E34').Value = '  +0.04%  '
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.603'
$ws.Range('E35').Value = '  -0.31%  '
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9939'
$ws.Range('E36').Value = '  -0.41%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.6255'
$ws.Range('E37').Value = '  -0.04%  '
$ws.Range('B38').Value = 'TrustWalletToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.9135'
$ws.Range('E38').Value = '  +0.14%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.392'
$ws.Range('E39').Value = '  +0.93%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.950'
$ws.Range('E40').Value = '  -5.92%  '
$ws.Range('B41').Value = 'PaxDollar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.001'
$ws.Range('E41').Value = '  +0.01%  '
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.01481'
$ws.Range('E42').Value = '  -1.47%  '
$ws.Range('B43').Value = 'Quant'
$ws.Range('C43').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '99.95'
$ws.Range('E43').Value = '  -3.26%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.436'
$ws.Range('E44').Value = '  -3.41%  '
$ws.Range('B45').Value = 'TheSandbox'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.3839'
$ws.Range('E45').Value = '  -0.69%  '
$ws.Range('B46').Value = 'Aptos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '6.693'
$ws.Range('E46').Value = '  -3.95%  '
$ws.Range('B47').Value = 'Algorand'
$ws.Range('C47').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.1156'
$ws.Range('E47').Value = '  -2.17%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.05365'
$ws.Range('E48').Value = '  +0.05%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '30.09'
$ws.Range('E49').Value = '  -1.28%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.689'
$ws.Range('E50').Value = '  -1.49%  '
$ws.Range('B51').Value = 'NEARProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.236'
$ws.Range('E51').Value = '  -1.47%  '
